$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" -------------------------------------------
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("P14").Value = 118.09
$wsVentasGrupo.Range("M25").Value = 4485.95
$wsVentasGrupo.Range("M35").Value = 13230.42
$wsVentasGrupo.Range("P55").Value = "1 de 53"

# --- Sheet 2: "VENTA MENSUAL" ----------------------------------------------
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F14").Value = 57
$wsVentaMensual.Range("F25").Value = 4485.95
$wsVentaMensual.Range("F35").Value = 13383.92
$wsVentaMensual.Range("F55").Value = 39056.7
# NOTE: this runtime's ColumnWidth -> OOXML <col width> conversion adds a
# fixed offset of 5/6 (0.8333333333333333) relative to the raw stored width,
# so we subtract it here to land on the exact target width of 14.
$targetColFWidth = 14 - (5 / 6)
$wsVentaMensual.Columns.Item(6).ColumnWidth = $targetColFWidth

# --- Sheet 3: "CUMPLIMIENTO MENSUAL" ---------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D10").Value = 118.09
$wsCumplimiento.Range("E10").Value = 1182.41
$wsCumplimiento.Range("F10").Value = 0.09080353710111495

$wsCumplimiento.Range("D16").Value = 34677.41
$wsCumplimiento.Range("E16").Value = 17149.05
$wsCumplimiento.Range("F16").Value = 0.6691062827752465

$wsCumplimiento.Range("D19").Value = 39056.7
$wsCumplimiento.Range("E19").Value = 74649.75064517916
$wsCumplimiento.Range("F19").Value = 0.3434871089405156
